# Berechnungen.xlsx edit script
# Implements:
#  - Foglio1 (sheet3): selection change J11 -> C6
#  - Foglio2 (sheet2): selection change M40 -> K31 (and drop topLeftCell),
#                       K30 formula made explicit
#  - Foglio3 (sheet4): new rows for "I Ball", "I Ball gesamt", "m rad",
#                       "I Rad", "I Gesamt", "alpha", "Drehmoment", K5 helper
#  - New sheet Foglio4 (sheet5): Drehzahl table with formatted header/table

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Foglio1: just move the selection
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Foglio1")
$ws1.Activate() | Out-Null
$ws1.Range("C6").Select() | Out-Null

# ---------------------------------------------------------------------------
# Foglio2: update K30 formula + move the selection
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Foglio2")
$ws2.Range("K30").Formula = "=(SQRT(C30)*I30)/(SQRT(2)*SQRT(I30*H30+A30-B30)*F30)"
$ws2.Activate() | Out-Null
$ws2.Range("K31").Select() | Out-Null

Write-Output "done"
